$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header: "email" -> "Email"
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "Email"

# ---------------------------------------------------------------------
# 2) Row 6: clear the email value (keeps the Hyperlink cell style, s=1,
#    but the cell itself becomes blank) and drop its hyperlink.
# ---------------------------------------------------------------------
$ws.Range("E6").ClearContents()

# ---------------------------------------------------------------------
# 3) Row 9: Nom/email edits
#    - C9 "Malika" -> long garbled string
# ---------------------------------------------------------------------
$ws.Range("C9").Value = "Malikakkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkk"

# ---------------------------------------------------------------------
# 4) New row 10 - plain values
# ---------------------------------------------------------------------
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "hhhhhh"
$ws.Range("C10").Value = "d"
$ws.Range("D10").Value = 645969744

# E10 / F10 need to carry over the same cell styles used by the other
# data rows (E: Hyperlink style s=1, F: short-date format s=2). Clone
# the style (format-only) from an existing cell via copy/paste-special
# through a scratch cell so we reuse the existing style slots instead
# of minting new ones.
$ws.Range("F2").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("ZZ1").Clear()
$ws.Range("F10").Value = 37582

# ---------------------------------------------------------------------
# 5) Hyperlinks - rebuild the list:
#    - drop the one for E6 (now blank)
#    - keep the rest pointing at the same mailto targets
#    - E9's target changes to "Malika#gmail.com" but keeps showing the
#      old "Malika@gmail.com" text (a stale "display" - as if someone
#      edited the link address without updating the visible text)
#    - add one for the new E10 cell
#    Hyperlinks.Add() overwrites the cell's displayed text with its
#    TextToDisplay argument (or the raw address, if omitted on a blank
#    cell), and always re-applies the built-in Hyperlink style to the
#    target cell (minting a redundant style slot) - so the real text is
#    re-asserted and the style re-stamped (s=1) afterwards.
# ---------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:Ahmed@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:Fatima@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:hajar@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:mohmed@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:gg@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:dfdd@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:Malika#gmail.com", "", "", "Malika@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E10"), "mailto:hajar@gmail.com")

# Re-stamp the Hyperlink cell style (s=1) on every touched E-cell.
$ws.Range("E2").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)
foreach ($addr in @("E2","E3","E4","E5","E7","E8","E9","E10")) {
    $ws.Range("ZZ1").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$ws.Range("ZZ1").Clear()

# Re-assert the real cell text (Hyperlinks.Add above clobbers it).
$ws.Range("E2").Value = "hajar@gmail.com"
$ws.Range("E3").Value = "Ahmed@gmail.com"
$ws.Range("E4").Value = "mohmed@gmail.com"
$ws.Range("E5").Value = "Fatima@gmail.com"
$ws.Range("E7").Value = "gg@gmail.com"
$ws.Range("E8").Value = "dfdd@gmail.com"
$ws.Range("E9").Value = "Malika#gmail.com"
$ws.Range("E10").Value = "hajar@gmail.com"
